# Daily attendance processing - 2026-01-10 01:38:19
#
# The "Recorded By" column (G) lists the users who recorded/touched each
# attendance session, as a comma-separated string. This pass normalizes
# the ordering of that list by rotating it left one position (the first
# name in the list is moved to the end), for the rows below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recordedBy = [ordered]@{
    2   = "backup@backdoor.com, system, System"
    3   = "dnasr281@gmail.com, System"
    4   = "backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    6   = "dnasr281@gmail.com, System"
    8   = "backup@backdoor.com, System"
    28  = "backup@backdoor.com, system, System"
    29  = "dnasr281@gmail.com, System"
    30  = "backup@backdoor.com, System"
    31  = "backup@backdoor.com, System"
    32  = "dnasr281@gmail.com, System"
    34  = "backup@backdoor.com, System"
    54  = "backup@backdoor.com, system, System"
    55  = "dnasr281@gmail.com, System"
    56  = "backup@backdoor.com, System"
    57  = "backup@backdoor.com, System"
    58  = "dnasr281@gmail.com, System"
    60  = "backup@backdoor.com, System"
    80  = "backup@backdoor.com, System"
    81  = "backup@backdoor.com, System"
    82  = "backup@backdoor.com, System"
    87  = "dnasr281@gmail.com, admin@admin.com"
    106 = "backup@backdoor.com, System"
    107 = "backup@backdoor.com, System"
    108 = "backup@backdoor.com, System"
    113 = "dnasr281@gmail.com, admin@admin.com"
    132 = "backup@backdoor.com, System"
    133 = "backup@backdoor.com, System"
    134 = "backup@backdoor.com, System"
    139 = "dnasr281@gmail.com, admin@admin.com"
}

foreach ($row in $recordedBy.Keys) {
    $ws.Cells.Item($row, 7).Value = $recordedBy[$row]
}
